# Add a new "Thank You" closing slide (slide 3) after the existing two
# slides, using the same "Title and Content" layout as slide 2.

$p = $ppt.ActivePresentation

# 2 = ppLayoutText / "Title and Content" custom layout (slideLayout2.xml) -
# matches the layout already used by slide 2.
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

# Give the two placeholders the same Korean auto-generated names PowerPoint
# assigns for a Korean-locale install ("제목 1" = Title 1, "내용 개체 틀 2" =
# Content Placeholder 2).
$title = $newSlide.Shapes.Item(1)
$content = $newSlide.Shapes.Item(2)
$title.Name = "제목 1"
$content.Name = "내용 개체 틀 2"

# Set the title text; the content placeholder is left empty, exactly like
# the authored slide.
$title.TextFrame.TextRange.Text = "Thank You"
